# DemoQATest.xlsx - "Test Steps" sheet keyword-driven-framework touch-up.
#
# The old "Enter Date" step (old row 9: "Step 8" / dateOfBirth) is dropped
# entirely (row deleted, rows below shift up one). A handful of cells in
# the remaining rows are corrected, and a new "Close Browser" step row is
# appended at the end so the sheet stays 16 rows tall.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Steps")

# 1. Remove the old "Enter Date" step row (old row 9) - everything below
#    shifts up by one row.
$ws.Range("A9").EntireRow.Delete()

# 2. Fill in the "Data" column for the Launch Browser / Navigate-to-url
#    steps, which used to be blank ("NA").
$ws.Range("H2").Value = "browser"
$ws.Range("H3").Value = "url"

# 3. "Select gender" is actually driven through a click action in the
#    keyword-driven framework now, not a select.
$ws.Range("E7").Value = "click"

# 4. The subjects locator value loses its trailing space (old row 10,
#    now row 9 after the delete above).
$ws.Range("G9").Value = "subjects"

# 5. "Select Hobbies" (old row 11, now row 10) is driven via click too,
#    and no longer carries a stray numeric data value.
$ws.Range("E10").Value = "click"
$ws.Range("H10").Value = "NA"

# 6. State/City steps (old rows 14/15, now 13/14 after the delete)
#    now carry real sample data instead of placeholder numbers.
$ws.Range("H13").Value = "NCR"
$ws.Range("H14").Value = "Delhi"

# 7. New final step: closing the browser.
$ws.Range("D16").Value = "Blose browser"
$ws.Range("E16").Value = "Close Browser"
$ws.Range("F16").Value = "NA"
$ws.Range("G16").Value = "NA"
$ws.Range("H16").Value = "NA"

# 8. Reflect the author's last on-screen selection position.
$ws.Range("G9").Select() | Out-Null
